$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.012.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.287.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.78%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.276.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.606"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.96%  "

$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "52.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.67%  "

$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.821.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.58%  "

$ws.Range("E17").Value = "  -4.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.292.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "62.989.86"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.959"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "414.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.94%  "

$ws.Range("E24").Value = "  -1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "569.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.55%  "

$ws.Range("E34").Value = "  -3.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.88%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.144"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.84%  "

$ws.Range("E38").Value = "  +4.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0734"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.95%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.102.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.44%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.360"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.73%  "

$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("E46").Value = "  -3.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.69%  "

$ws.Range("E48").Value = "  -3.52%  "

$ws.Range("E49").Value = "  -5.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.43%  "
